$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scraped-data refresh: updated Price (D) / Volume(1h) (E) figures and
# the Hora (G) timestamp (14 -> 15) for each coin row. Only cells whose
# text actually changes are touched.
$changes = @(
    ,@("D2", "281.75")
    ,@("E2", "1.26%")
    ,@("G2", "15")
    ,@("D3", "28.19")
    ,@("E3", "3.02%")
    ,@("G3", "15")
    ,@("D4", "5.051")
    ,@("E4", "4.13%")
    ,@("G4", "15")
    ,@("D5", "0.06488")
    ,@("E5", "1.71%")
    ,@("G5", "15")
    ,@("D6", "7.224")
    ,@("E6", "2.99%")
    ,@("G6", "15")
    ,@("D7", "3.379")
    ,@("E7", "2.32%")
    ,@("G7", "15")
    ,@("D8", "1.391")
    ,@("E8", "3.03%")
    ,@("G8", "15")
    ,@("D9", "0.9281")
    ,@("E9", "5.85%")
    ,@("G9", "15")
    ,@("D10", "0.1553")
    ,@("E10", "2.31%")
    ,@("G10", "15")
    ,@("D11", "0.06255")
    ,@("E11", "19.32%")
    ,@("G11", "15")
    ,@("D12", "0.07543")
    ,@("E12", "0.78%")
    ,@("G12", "15")
    ,@("D13", "0.02858")
    ,@("E13", "-1.81%")
    ,@("G13", "15")
    ,@("D14", "0.08972")
    ,@("E14", "0.18%")
    ,@("G14", "15")
    ,@("D15", "0.001579")
    ,@("E15", "1.17%")
    ,@("G15", "15")
    ,@("D16", "0.0006393")
    ,@("E16", "0.40%")
    ,@("G16", "15")
    ,@("D17", "0.006163")
    ,@("E17", "1.09%")
    ,@("G17", "15")
    ,@("E18", "-0.94%")
    ,@("G18", "15")
    ,@("D19", "2.230")
    ,@("E19", "-0.69%")
    ,@("G19", "15")
    ,@("D20", "0.3189")
    ,@("E20", "3.12%")
    ,@("G20", "15")
    ,@("D21", "0.1279")
    ,@("E21", "-3.95%")
    ,@("G21", "15")
    ,@("D22", "4.039")
    ,@("E22", "3.42%")
    ,@("G22", "15")
    ,@("D23", "0.1544")
    ,@("E23", "0.55%")
    ,@("G23", "15")
    ,@("D24", "0.04390")
    ,@("E24", "-0.16%")
    ,@("G24", "15")
    ,@("D25", "0.001183")
    ,@("E25", "1.10%")
    ,@("G25", "15")
    ,@("D26", "0.004396")
    ,@("E26", "12.81%")
    ,@("G26", "15")
    ,@("D27", "0.0001250")
    ,@("E27", "5.97%")
    ,@("G27", "15")
    ,@("D28", "0.0001618")
    ,@("E28", "-1.82%")
    ,@("G28", "15")
    ,@("G29", "15")
    ,@("G30", "15")
    ,@("G31", "15")
    ,@("G32", "15")
    ,@("G33", "15")
    ,@("G34", "15")
    ,@("G35", "15")
    ,@("G36", "15")
    ,@("G37", "15")
    ,@("G38", "15")
    ,@("G39", "15")
    ,@("D40", "0.04149")
    ,@("E40", "1.93%")
    ,@("G40", "15")
    ,@("D41", "0.006674")
    ,@("E41", "-2.08%")
    ,@("G41", "15")
    ,@("D42", "0.1222")
    ,@("E42", "-13.64%")
    ,@("G42", "15")
    ,@("D43", "0.002051")
    ,@("E43", "8.53%")
    ,@("G43", "15")
    ,@("D44", "0.01208")
    ,@("E44", "3.42%")
    ,@("G44", "15")
    ,@("D45", "0.00005600")
    ,@("E45", "4.51%")
    ,@("G45", "15")
    ,@("G46", "15")
    ,@("D47", "0.01300")
    ,@("E47", "-29.75%")
    ,@("G47", "15")
    ,@("G48", "15")
    ,@("G49", "15")
    ,@("G50", "15")
    ,@("G51", "15")
)

foreach ($change in $changes) {
    $addr = $change[0]
    $val = $change[1]

    $cell = $ws.Range($addr)
    $origFormat = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = $origFormat
}
